# Updating location lat and long name, and fixing a habitat-type typo
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("attribute")

# Fix typo: "borad" -> "broad"
$ws.Range("B4").Value = "sample site habitat type, broad"

# Fix the attribute name for lat/lon
$ws.Range("A3").Value = "Lat Long UTM"

# Update the selected cell on this sheet
$ws.Range("C20").Select()

# Extend existing data validations to include row 1, and add two new validations
$ws.Range("C1:C6").Validation.Delete()
$ws.Range("C1:C6").Validation.Add(3, 1, 1, """string,boolean,decimal,float,double,duration,dateTime,time,date,gYearMonth,gYear,gMonthDay,gDay,gMonth""")

$ws.Range("D1:D6").Validation.Delete()
$ws.Range("D1:D6").Validation.Add(3, 1, 1, """nominal,ordinal,interval,ratio,dateTime""")

$ws.Range("F1:F6").Validation.Delete()
$ws.Range("F1:F6").Validation.Add(3, 1, 1, """text,enumerated,dateTime,numeric""")

$ws.Range("G1").Validation.Add(3, 1, 1, """ratio,interval""")

$ws.Range("I1").Validation.Add(3, 1, 1, """natural,whole,interger,real""")
